# Fruta / hortaliza, semanal
# Apply a row-wise permutation to the data rows (2..39) across the columns
# that vary between records: D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg). Every other column is identical across all
# rows so this fully captures the edit described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# newRow -> oldRow : row newRow's new data comes from row oldRow's original data
$map = @{
    2  = 11
    3  = 6
    4  = 32
    5  = 39
    6  = 13
    7  = 37
    8  = 8
    9  = 28
    10 = 19
    11 = 26
    12 = 29
    13 = 5
    14 = 7
    15 = 23
    16 = 34
    17 = 33
    18 = 36
    19 = 10
    20 = 17
    21 = 12
    22 = 20
    23 = 15
    24 = 2
    25 = 25
    26 = 9
    27 = 21
    28 = 35
    29 = 18
    30 = 38
    31 = 14
    32 = 24
    33 = 31
    34 = 27
    35 = 16
    36 = 3
    37 = 4
    38 = 30
    39 = 22
}

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the original values for every touched column/row before writing
# anything, since several rows read from each other (it's a permutation).
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 39; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    if ($oldRow -eq $newRow) {
        continue
    }
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $snapshot["$col$oldRow"]
    }
}
